$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column. This shifts the existing Code / Description /
# Definition columns from A:C to B:D, matching the diff's D1/D2/... cells.
$ws.Columns("A").Insert()

# New header for the inserted column.
$ws.Range("A1").Value = "Version"

# Every data row gets a "1.0" version marker in the new column. Typing the
# literal string "1.0" via .Value would be auto-coerced to the number 1 (the
# same way real Excel treats unquoted numeric-looking input), so instead we
# build the text "1.0" with a formula, copy it, and paste-special just the
# values into each destination cell. That keeps the result a genuine text
# cell (t="s") without a quote-prefix cell style, matching the target.
$ws.Range("H1").Formula = "=""1.0"""
$ws.Range("H1").Copy()
$ws.Range("A2").PasteSpecial(-4163)
$ws.Range("A3").PasteSpecial(-4163)
$ws.Range("A4").PasteSpecial(-4163)
$ws.Range("A5").PasteSpecial(-4163)
$ws.Range("A6").PasteSpecial(-4163)
$ws.Range("A7").PasteSpecial(-4163)

# Remove the scratch helper cell used to build the text value.
$ws.Range("H1").Clear()
